$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45190 = 2023-09-21)
# that was updated to 45192 (2023-09-23) for every data row (2..266).
$newDate = 45192

for ($r = 2; $r -le 266; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
